# Apply the edits described by the diff:
# - Sheet "output_2": a few numeric result values changed slightly
# - Sheets "output_5", "output_6", "output_7", "output_8": the regression
#   summary's Date/Time stamp was refreshed to a new run timestamp

$wb = $excel.ActiveWorkbook

# --- output_2 sheet: update numeric results ---
$wsResults = $wb.Worksheets.Item("output_2")
$wsResults.Range("B15").Value = 8502247
$wsResults.Range("D17").Value = 6988263.5
$wsResults.Range("D18").Value = 7997392.5
$wsResults.Range("D19").Value = 10388613

# --- output_5/output_6/output_7/output_8: refresh Date/Time stamps ---
$sheetNames = @("output_5", "output_6", "output_7", "output_8")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B8").Value = "Tue, 07 Mar 2023"
    $ws.Range("B9").Value = "10:45:02"
}
